$ws = $excel.ActiveWorkbook.ActiveSheet
$ws.Range("A21").Value = "Batken oblast"
"ok"
